$wb = $excel.ActiveWorkbook

# --- New item identifiers / constants -------------------------------------------------
$md1        = "3d5354a3-5fdf-4fce-b931-6b41a3ccc432.md"
$md1Disp    = "e2e\3d5354a3-5fdf-4fce-b931-6b41a3ccc432.md"
$md2        = "50c7c997-a27f-46c0-a0fc-f5ef91ff8bed.md"
$md2Disp    = "e2e\50c7c997-a27f-46c0-a0fc-f5ef91ff8bed.md"

$zhXlf1     = "3d5354a3-5fdf-4fce-b931-6b41a3ccc432.ce536d311933389fb9de397dfa3926e172c73198.zh-cn.xlf"
$zhXlf2     = "50c7c997-a27f-46c0-a0fc-f5ef91ff8bed.698ad3621bd04f6c1388c20378aec8a84bad2388.zh-cn.xlf"
$deXlf1     = "3d5354a3-5fdf-4fce-b931-6b41a3ccc432.ce536d311933389fb9de397dfa3926e172c73198.de-de.xlf"
$deXlf2     = "50c7c997-a27f-46c0-a0fc-f5ef91ff8bed.698ad3621bd04f6c1388c20378aec8a84bad2388.de-de.xlf"

$hoDate     = "2016-08-28 22:41:00"
$zhGenDate  = "2016-08-28 22:40:55"
$status     = "Ready for handoff"
$ext        = ".md"
$epoch      = "0001-01-01 00:00:00"

$url1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/HEAD/e2e/$md1"
$url2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/HEAD/e2e/$md2"

# ===========================================================================
# Sheet "Overview" -- two new summary rows (4 and 5)
# ===========================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(4,1).Value = $md1
$wsOverview.Cells.Item(4,2).Value = $md1Disp
$wsOverview.Cells.Item(4,3).Value = $ext
$wsOverview.Cells.Item(4,4).Value = ""
$wsOverview.Cells.Item(4,5).Value = $status
$wsOverview.Cells.Item(4,6).Value = $status
$wsOverview.Cells.Item(4,7).Value = $hoDate

$wsOverview.Cells.Item(5,1).Value = $md2
$wsOverview.Cells.Item(5,2).Value = $md2Disp
$wsOverview.Cells.Item(5,3).Value = $ext
$wsOverview.Cells.Item(5,4).Value = ""
$wsOverview.Cells.Item(5,5).Value = $status
$wsOverview.Cells.Item(5,6).Value = $status
$wsOverview.Cells.Item(5,7).Value = $hoDate

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), $url1, $null, $null, $md1Disp) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), $url2, $null, $null, $md2Disp) | Out-Null

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G5"))

# ===========================================================================
# Sheet "zh-cn" -- two new detail rows (4 and 5)
# ===========================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Cells.Item(4,1).Value  = $md1
$wsZh.Cells.Item(4,2).Value  = $ext
$wsZh.Cells.Item(4,3).Value  = $status
$wsZh.Cells.Item(4,4).Value  = "e2e"
$wsZh.Cells.Item(4,5).Value  = "ht"
$wsZh.Cells.Item(4,6).Value  = "False"
$wsZh.Cells.Item(4,7).Value  = $zhXlf1
$wsZh.Cells.Item(4,8).Value  = $zhGenDate
$wsZh.Cells.Item(4,9).Value  = ""
$wsZh.Cells.Item(4,10).Value = ""
$wsZh.Cells.Item(4,11).Value = $epoch
$wsZh.Cells.Item(4,12).Value = ""
$wsZh.Cells.Item(4,13).Value = "True"
$wsZh.Cells.Item(4,14).Value = ""
$wsZh.Cells.Item(4,15).Value = "False"
$wsZh.Cells.Item(4,16).Value = ""

$wsZh.Cells.Item(5,1).Value  = $md2
$wsZh.Cells.Item(5,2).Value  = $ext
$wsZh.Cells.Item(5,3).Value  = $status
$wsZh.Cells.Item(5,4).Value  = "e2e"
$wsZh.Cells.Item(5,5).Value  = "ht"
$wsZh.Cells.Item(5,6).Value  = "False"
$wsZh.Cells.Item(5,7).Value  = $zhXlf2
$wsZh.Cells.Item(5,8).Value  = $zhGenDate
$wsZh.Cells.Item(5,9).Value  = ""
$wsZh.Cells.Item(5,10).Value = ""
$wsZh.Cells.Item(5,11).Value = $epoch
$wsZh.Cells.Item(5,12).Value = ""
$wsZh.Cells.Item(5,13).Value = "True"
$wsZh.Cells.Item(5,14).Value = ""
$wsZh.Cells.Item(5,15).Value = "False"
$wsZh.Cells.Item(5,16).Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $url1, $null, $null, $md1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), $url2, $null, $null, $md2) | Out-Null

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P5"))

# ===========================================================================
# Sheet "de-de" -- two new detail rows (4 and 5)
# ===========================================================================
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Cells.Item(4,1).Value  = $md1
$wsDe.Cells.Item(4,2).Value  = $ext
$wsDe.Cells.Item(4,3).Value  = $status
$wsDe.Cells.Item(4,4).Value  = "e2e"
$wsDe.Cells.Item(4,5).Value  = "ht"
$wsDe.Cells.Item(4,6).Value  = "False"
$wsDe.Cells.Item(4,7).Value  = $deXlf1
$wsDe.Cells.Item(4,8).Value  = $hoDate
$wsDe.Cells.Item(4,9).Value  = ""
$wsDe.Cells.Item(4,10).Value = ""
$wsDe.Cells.Item(4,11).Value = $epoch
$wsDe.Cells.Item(4,12).Value = ""
$wsDe.Cells.Item(4,13).Value = "True"
$wsDe.Cells.Item(4,14).Value = ""
$wsDe.Cells.Item(4,15).Value = "False"
$wsDe.Cells.Item(4,16).Value = ""

$wsDe.Cells.Item(5,1).Value  = $md2
$wsDe.Cells.Item(5,2).Value  = $ext
$wsDe.Cells.Item(5,3).Value  = $status
$wsDe.Cells.Item(5,4).Value  = "e2e"
$wsDe.Cells.Item(5,5).Value  = "ht"
$wsDe.Cells.Item(5,6).Value  = "False"
$wsDe.Cells.Item(5,7).Value  = $deXlf2
$wsDe.Cells.Item(5,8).Value  = $hoDate
$wsDe.Cells.Item(5,9).Value  = ""
$wsDe.Cells.Item(5,10).Value = ""
$wsDe.Cells.Item(5,11).Value = $epoch
$wsDe.Cells.Item(5,12).Value = ""
$wsDe.Cells.Item(5,13).Value = "True"
$wsDe.Cells.Item(5,14).Value = ""
$wsDe.Cells.Item(5,15).Value = "False"
$wsDe.Cells.Item(5,16).Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $url1, $null, $null, $md1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), $url2, $null, $null, $md2) | Out-Null

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P5"))

Write-Host "Generate Report for Handoff: added 2 rows to each sheet"
